$wb = $excel.ActiveWorkbook

# Rename the "T_Shirts" worksheet to "T-Shirts"
$ws = $wb.Worksheets.Item("T_Shirts")
$ws.Name = "T-Shirts"
